$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (font/border/alignment) from the existing header cell G1
# onto the new header cell H1, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Fill in the new "Save" data column values for the two data rows.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
